# Update "想去人数" (F column) figures across the workbook's sheets
# to reflect the refreshed scrape at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 520
$ws1.Range("F3").Value = 743
$ws1.Range("F4").Value = 1478
$ws1.Range("F5").Value = 223
$ws1.Range("F7").Value = 141
$ws1.Range("F8").Value = 6196
$ws1.Range("F12").Value = 5099
$ws1.Range("F15").Value = 1175
$ws1.Range("F16").Value = 55
$ws1.Range("F18").Value = 63
$ws1.Range("F20").Value = 293
$ws1.Range("F22").Value = 3612

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 77

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 77
$ws4.Range("F3").Value = 520
$ws4.Range("F4").Value = 743
$ws4.Range("F5").Value = 1478
$ws4.Range("F6").Value = 223
$ws4.Range("F8").Value = 141
$ws4.Range("F9").Value = 6196
$ws4.Range("F13").Value = 5099
$ws4.Range("F16").Value = 1175
$ws4.Range("F17").Value = 55
$ws4.Range("F19").Value = 63
$ws4.Range("F21").Value = 293
$ws4.Range("F23").Value = 3612
